$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

    $ws.Range("B2").Value = 1.02
    $ws.Range("C2").Value = 1.033400520885426
    $ws.Range("D2").Value = 1.034327182574525
    $ws.Range("E2").Value = 1.04148653396911
    $ws.Range("F2").Value = 1.049490449205274
    $ws.Range("I2").Value = 1.033290551877651
    $ws.Range("J2").Value = 1.038525453572313
    $ws.Range("K2").Value = 1.037126975327753
    $ws.Range("L2").Value = 1.044265911982581
    $ws.Range("M2").Value = 1.052247361867015
    $ws.Range("N2").Value = 1.04000027883692
    $ws.Range("B3").Value = 1.02
    $ws.Range("C3").Value = 1.034408020088846
    $ws.Range("D3").Value = 1.035198260687544
    $ws.Range("E3").Value = 1.042416546466402
    $ws.Range("F3").Value = 1.050587514690108
    $ws.Range("I3").Value = 1.033441152619347
    $ws.Range("J3").Value = 1.039175356618244
    $ws.Range("K3").Value = 1.037807053238524
    $ws.Range("L3").Value = 1.045006248652858
    $ws.Range("M3").Value = 1.053155948611114
    $ws.Range("N3").Value = 1.040651104819721
    $ws.Range("B4").Value = 1.02
    $ws.Range("C4").Value = 1.035059673641464
    $ws.Range("D4").Value = 1.035761960277996
    $ws.Range("E4").Value = 1.043018542445109
    $ws.Range("F4").Value = 1.051297902548048
    $ws.Range("I4").Value = 1.033536319260591
    $ws.Range("J4").Value = 1.039595039987857
    $ws.Range("K4").Value = 1.038246527103056
    $ws.Range("L4").Value = 1.045484886589158
    $ws.Range("M4").Value = 1.053743786209932
    $ws.Range("N4").Value = 1.041071384187858
    $ws.Range("B5").Value = 1.02
    $ws.Range("C5").Value = 1.035333564906685
    $ws.Range("D5").Value = 1.035998951683263
    $ws.Range("E5").Value = 1.043271672516164
    $ws.Range("F5").Value = 1.051596671969141
    $ws.Range("I5").Value = 1.033575780349301
    $ws.Range("J5").Value = 1.039771271441099
    $ws.Range("K5").Value = 1.038431141964318
    $ws.Range("L5").Value = 1.045686007411659
    $ws.Range("M5").Value = 1.053990893762609
    $ws.Range("N5").Value = 1.041247865909977
    $ws.Range("B6").Value = 1.02
    $ws.Range("C6").Value = 1.035379548701312
    $ws.Range("D6").Value = 1.036038744269469
    $ws.Range("E6").Value = 1.043314177137457
    $ws.Range("F6").Value = 1.051646843831131
    $ws.Range("I6").Value = 1.033582373946906
    $ws.Range("J6").Value = 1.039800849531084
    $ws.Range("K6").Value = 1.038462131387007
    $ws.Range("L6").Value = 1.04571977068646
    $ws.Range("M6").Value = 1.054032383056864
    $ws.Range("N6").Value = 1.041277486004243
    $ws.Range("B7").Value = 1.02
    $ws.Range("C7").Value = 1.035063333641328
    $ws.Range("D7").Value = 1.035765126922221
    $ws.Range("E7").Value = 1.043021924582132
    $ws.Range("F7").Value = 1.051301894239479
    $ws.Range("I7").Value = 1.033536848691605
    $ws.Range("J7").Value = 1.039597395601267
    $ws.Range("K7").Value = 1.038248994486574
    $ws.Range("L7").Value = 1.045487574363446
    $ws.Range("M7").Value = 1.053747088149934
    $ws.Range("N7").Value = 1.041073743146509
    $ws.Range("B8").Value = 1.02
    $ws.Range("C8").Value = 1.03374106484379
    $ws.Range("D8").Value = 1.034621555892857
    $ws.Range("E8").Value = 1.041800790978546
    $ws.Range("F8").Value = 1.049861101089732
    $ws.Range("I8").Value = 1.033341920256333
    $ws.Range("J8").Value = 1.038745266457644
    $ws.Range("K8").Value = 1.037356930986089
    $ws.Range("L8").Value = 1.044516196662694
    $ws.Range("M8").Value = 1.052554439016831
    $ws.Range("N8").Value = 1.040220403881761
    $ws.Range("B9").Value = 1.02
    $ws.Range("C9").Value = 1.031409032854192
    $ws.Range("D9").Value = 1.032606883502903
    $ws.Range("E9").Value = 1.039650678779623
    $ws.Range("F9").Value = 1.047326189755839
    $ws.Range("I9").Value = 1.032980980173104
    $ws.Range("J9").Value = 1.037237240908386
    $ws.Range("K9").Value = 1.035780569938162
    $ws.Range("L9").Value = 1.042801393383744
    $ws.Range("M9").Value = 1.05045225738538
    $ws.Range("N9").Value = 1.03871023676325
    $ws.Range("B10").Value = 1.02
    $ws.Range("C10").Value = 1.029852993312997
    $ws.Range("D10").Value = 1.031264111353209
    $ws.Range("E10").Value = 1.038218437644567
    $ws.Range("F10").Value = 1.045638928930262
    $ws.Range("I10").Value = 1.032728655019835
    $ws.Range("J10").Value = 1.036227574478175
    $ws.Range("K10").Value = 1.034726716475923
    $ws.Range("L10").Value = 1.041656128682834
    $ws.Range("M10").Value = 1.049050438586347
    $ws.Range("N10").Value = 1.037699136490899
    $ws.Range("B11").Value = 1.02
    $ws.Range("C11").Value = 1.029178891435843
    $ws.Range("D11").Value = 1.030682765748276
    $ws.Range("E11").Value = 1.037598545083382
    $ws.Range("F11").Value = 1.044908966088334
    $ws.Range("I11").Value = 1.032616627585842
    $ws.Range("J11").Value = 1.03578935974163
    $ws.Range("K11").Value = 1.034269693755286
    $ws.Range("L11").Value = 1.041159732664638
    $ws.Range("M11").Value = 1.048443353817696
    $ws.Range("N11").Value = 1.037260299439163
    $ws.Range("B12").Value = 1.02
    $ws.Range("C12").Value = 1.028928450878514
    $ws.Range("D12").Value = 1.030466841243492
    $ws.Range("E12").Value = 1.037368331507013
    $ws.Range("F12").Value = 1.044637920586188
    $ws.Range("I12").Value = 1.032574600081634
    $ws.Range("J12").Value = 1.035626434020272
    $ws.Range("K12").Value = 1.034099830833151
    $ws.Range("L12").Value = 1.040975275984306
    $ws.Range("M12").Value = 1.048217842557839
    $ws.Range("N12").Value = 1.037097142344593
    $ws.Range("B13").Value = 1.02
    $ws.Range("C13").Value = 1.0289821734513
    $ws.Range("D13").Value = 1.030513157179074
    $ws.Range("E13").Value = 1.037417711187233
    $ws.Range("F13").Value = 1.044696056455178
    $ws.Range("I13").Value = 1.032583633919826
    $ws.Range("J13").Value = 1.035661389070011
    $ws.Range("K13").Value = 1.034136271721372
    $ws.Range("L13").Value = 1.041014845870714
    $ws.Range("M13").Value = 1.048266216070166
    $ws.Range("N13").Value = 1.037132147034512
    $ws.Range("B14").Value = 1.02
    $ws.Range("C14").Value = 1.029158190943334
    $ws.Range("D14").Value = 1.030664917090049
    $ws.Range("E14").Value = 1.037579514695446
    $ws.Range("F14").Value = 1.044886559439335
    $ws.Range("I14").Value = 1.032613162053121
    $ws.Range("J14").Value = 1.035775895374997
    $ws.Range("K14").Value = 1.034255654959187
    $ws.Range("L14").Value = 1.041144486907086
    $ws.Range("M14").Value = 1.048424713249168
    $ws.Range("N14").Value = 1.037246815951586
    $ws.Range("B15").Value = 1.02
    $ws.Range("C15").Value = 1.029266634620734
    $ws.Range("D15").Value = 1.03075842313785
    $ws.Range("E15").Value = 1.037679212772185
    $ws.Range("F15").Value = 1.045003947236319
    $ws.Range("I15").Value = 1.032631300263405
    $ws.Range("J15").Value = 1.035846426191409
    $ws.Range("K15").Value = 1.034329197099021
    $ws.Range("L15").Value = 1.041224353352432
    $ws.Range("M15").Value = 1.048522366884692
    $ws.Range("N15").Value = 1.037317446929847
    $ws.Range("B16").Value = 1.02
    $ws.Range("C16").Value = 1.029897724351269
    $ws.Range("D16").Value = 1.031302695157222
    $ws.Range("E16").Value = 1.03825958375418
    $ws.Range("F16").Value = 1.045687387498751
    $ws.Range("I16").Value = 1.032736031613735
    $ws.Range("J16").Value = 1.036256635835178
    $ws.Range("K16").Value = 1.034757032914722
    $ws.Range("L16").Value = 1.041689062582058
    $ws.Range("M16").Value = 1.049090726972798
    $ws.Range("N16").Value = 1.037728239118362
    $ws.Range("B17").Value = 1.02
    $ws.Range("C17").Value = 1.03029350282014
    $ws.Range("D17").Value = 1.031644125407862
    $ws.Range("E17").Value = 1.038623709707166
    $ws.Range("F17").Value = 1.046116260938226
    $ws.Range("I17").Value = 1.032800985831324
    $ws.Range("J17").Value = 1.036513676012214
    $ws.Range("K17").Value = 1.035025216649239
    $ws.Range("L17").Value = 1.041980431773811
    $ws.Range("M17").Value = 1.049447220880161
    $ws.Range("N17").Value = 1.037985644321931
    $ws.Range("B18").Value = 1.02
    $ws.Range("C18").Value = 1.03052432220781
    $ws.Range("D18").Value = 1.031843283924211
    $ws.Range("E18").Value = 1.038836124940244
    $ws.Range("F18").Value = 1.046366476546617
    $ws.Range("I18").Value = 1.032838605329689
    $ws.Range("J18").Value = 1.036663504517503
    $ws.Range("K18").Value = 1.035181576385669
    $ws.Range("L18").Value = 1.042150335304562
    $ws.Range("M18").Value = 1.049655149238567
    $ws.Range("N18").Value = 1.038135685600883
    $ws.Range("B19").Value = 1.02
    $ws.Range("C19").Value = 1.030603020313171
    $ws.Range("D19").Value = 1.031911193175475
    $ws.Range("E19").Value = 1.038908557564352
    $ws.Range("F19").Value = 1.046451803940788
    $ws.Range("I19").Value = 1.032851387270303
    $ws.Range("J19").Value = 1.036714575416636
    $ws.Range("K19").Value = 1.035234879587342
    $ws.Range("L19").Value = 1.042208260003076
    $ws.Range("M19").Value = 1.049726045984231
    $ws.Range("N19").Value = 1.038186829026551
    $ws.Range("B20").Value = 1.02
    $ws.Range("C20").Value = 1.030251042794172
    $ws.Range("D20").Value = 1.0316074923386
    $ws.Range("E20").Value = 1.038584639670703
    $ws.Range("F20").Value = 1.046070240567116
    $ws.Range("I20").Value = 1.032794044492477
    $ws.Range("J20").Value = 1.036486108238137
    $ws.Range("K20").Value = 1.034996450030498
    $ws.Range("L20").Value = 1.041949175495413
    $ws.Range("M20").Value = 1.049408973335485
    $ws.Range("N20").Value = 1.037958037398454
    $ws.Range("B21").Value = 1.02
    $ws.Range("C21").Value = 1.029106359535085
    $ws.Range("D21").Value = 1.030620227211535
    $ws.Range("E21").Value = 1.037531866428706
    $ws.Range("F21").Value = 1.044830458429023
    $ws.Range("I21").Value = 1.032604478220611
    $ws.Range("J21").Value = 1.03574218034993
    $ws.Range("K21").Value = 1.034220502448065
    $ws.Range("L21").Value = 1.041106312868187
    $ws.Range("M21").Value = 1.048378040135628
    $ws.Range("N21").Value = 1.037213053047316
    $ws.Range("B22").Value = 1.02
    $ws.Range("C22").Value = 1.028386367594595
    $ws.Range("D22").Value = 1.029999570729357
    $ws.Range("E22").Value = 1.036870190031362
    $ws.Range("F22").Value = 1.044051509084536
    $ws.Range("I22").Value = 1.032482886500001
    $ws.Range("J22").Value = 1.035273556837834
    $ws.Range("K22").Value = 1.033732029955788
    $ws.Range("L22").Value = 1.040575948929649
    $ws.Range("M22").Value = 1.047729776740822
    $ws.Range("N22").Value = 1.03674376403608
    $ws.Range("B23").Value = 1.02
    $ws.Range("C23").Value = 1.028768075898721
    $ws.Range("D23").Value = 1.030328585052164
    $ws.Range("E23").Value = 1.037220933936271
    $ws.Range("F23").Value = 1.044464392582574
    $ws.Range("I23").Value = 1.032547572274136
    $ws.Range("J23").Value = 1.03552206694128
    $ws.Range("K23").Value = 1.033991035530264
    $ws.Range("L23").Value = 1.040857144846731
    $ws.Range("M23").Value = 1.048073440431106
    $ws.Range("N23").Value = 1.036992627052378
    $ws.Range("B24").Value = 1.02
    $ws.Range("C24").Value = 1.030270228749379
    $ws.Range("D24").Value = 1.031624045220899
    $ws.Range("E24").Value = 1.038602293655788
    $ws.Range("F24").Value = 1.046091035003331
    $ws.Range("I24").Value = 1.032797181810303
    $ws.Range("J24").Value = 1.036498565233437
    $ws.Range("K24").Value = 1.035009448635143
    $ws.Range("L24").Value = 1.041963299007941
    $ws.Range("M24").Value = 1.049426255780755
    $ws.Range("N24").Value = 1.037970512084116
    $ws.Range("B25").Value = 1.02
    $ws.Range("C25").Value = 1.032012157957061
    $ws.Range("D25").Value = 1.033127667299046
    $ws.Range("E25").Value = 1.040206331277004
    $ws.Range("F25").Value = 1.047981054452819
    $ws.Range("I25").Value = 1.033076356272996
    $ws.Range("J25").Value = 1.037627864606436
    $ws.Range("K25").Value = 1.036188618716203
    $ws.Range("L25").Value = 1.039101415191756
